$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Juan"
$ws.Range("B2").Value = "Alberto"
$ws.Range("C2").Value = "Pérez"
$ws.Range("D2").Value = "Cáceres"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "20230001"
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = "Introducción a la Informática"
$ws.Range("H2").Value = 800

$ws.Range("A3").Value = "David"
$ws.Range("B3").Value = "José"
$ws.Range("C3").Value = "Martínez"
$ws.Range("D3").Value = "Rodríguez"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "201905678"
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = "Introducción a la Informática"
$ws.Range("H3").Value = 800
